$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the "ECs" target-cluster row (original row 2); remaining rows shift up.
$ws.Rows.Item(2).Delete()

# Row 2 (was row 3, Inflammatory-Mac) - update recalculated values
$ws.Range("G2").Value = 0.8243956666666667
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 0.6666666666666666
$ws.Range("M2").Value = 1.375733333333333
$ws.Range("N2").Value = 4.1272
$ws.Range("O2").Value = 0.457732955319909
$ws.Range("P2").Value = 0.457732955319909
$ws.Range("Q2").Value = 1.134148598488889
$ws.Range("R2").Value = 10.2073373864
$ws.Range("S2").Value = 0.457732955319909
$ws.Range("T2").Value = 0.457732955319909

# Row 3 (was row 4, Resolving-Mac) - update recalculated values
$ws.Range("G3").Value = 0.8243956666666667
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 1.629803666666667
$ws.Range("N3").Value = 4.889411000000001
$ws.Range("O3").Value = 0.542267044680091
$ws.Range("P3").Value = 0.542267044680091
$ws.Range("Q3").Value = 1.343603080317445
$ws.Range("R3").Value = 12.092427722857
$ws.Range("S3").Value = 0.542267044680091
$ws.Range("T3").Value = 0.542267044680091
